$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 (use same style as other headers via Style copy)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold, bordered, centered) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill I2:I13 with 1, and J2:J13 with the same values as H2:H13
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value = $hVal
}
